$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some columns (L, M) carry a Text ("@") number format, so a plain
# `.Value = 0` assignment gets stored as the literal string "0" instead of
# a real number. Round-trip the format through General while writing the
# value, then restore the original format, so the written value stays a
# genuine number (matching how the rest of the column is stored) without
# permanently changing the cell's display format.
function Set-NumericValue($sheet, $addr, $val) {
    $range = $sheet.Range($addr)
    $fmt = $range.NumberFormat
    $range.NumberFormat = "General"
    $range.Value = $val
    $range.NumberFormat = $fmt
}

# --- Daily data updates (new counts arrived for already-entered days) ---
# Row 310 (2020-xx-xx): new positive cases revised 130 -> 131
$ws.Range("C310").Value = 131

# Row 316: new positive cases revised 144 -> 145
$ws.Range("C316").Value = 145

# Row 322: new positive cases revised 95 -> 149; one new extra-hospital death recorded
$ws.Range("C322").Value = 149
Set-NumericValue $ws "M322" 1

# Row 323: new positive cases revised 21 -> 105; one new extra-hospital death recorded
$ws.Range("C323").Value = 105
Set-NumericValue $ws "M323" 1

# --- Row 324: fill in the day's figures (previously blank placeholder row) ---
$ws.Range("C324").Value = 16
$ws.Range("E324").Value = 10
$ws.Range("F324").Value = 8
$ws.Range("G324").Value = 114

# L324/M324 are formatted as Text ("@"); see Set-NumericValue above.
Set-NumericValue $ws "L324" 0
Set-NumericValue $ws "M324" 0

# --- Keep the last selected cell on the data sheet pointed at column Z ---
$ws.Range("Z2").Select()
